{"js": "// Office.js (Word JavaScript API) script\n//\n// The source document has a single paragraph containing one long run of text.\n// The target rewrites that paragraph as the opening paragraph of a new (much\n// shorter) article, then appends six more paragraphs (one per quoted clause/\n// source grouping) and a final empty paragraph.\n\nconst paragraphTexts = [\n  \"After a frightening study from last week showed that industrial methane emissions have been \\\"vastly underestimated,\\\" a new projection Friday that the United States is on track to become the world's leading exporter of liquefied natural gas within five years provoked warnings that the American fracking boom could \\\"end hope for climate stability.\\\" \",\n  \"Liquefied natural gas (LNG) is primarily composed of methane, a greenhouse gas that is 84 times more potent than carbon dioxide over a 20-year period. Methane emissions, by some estimates, are responsible for about a quarter of human-caused global warming.  \\\"Science confirms that gas is a climate killer,\\\" Wenonah Hauter, executive director of the U.S. advocacy group Food & Water Watch, said in a statement Friday, citing methane's planet-warming potential. \",\n  \"Hauter's statement came in response to an International Energy Agency (IEA) annual report, released Friday, that featured the new projection about U.S. LNG exports. The IEA report states that global demand for natural gas grew last year at the fastest rate in nearly a decade and is expected to keep growing, \\\"driven by strong consumption in fast-growing Asian economies and supported by the continued development of the international gas trade.\\\"  The IEA's new release came just two days after Food & Water Watch published a report which, as Hauter put it, \\\"shows that the power, petrochemical, and LNG export industries are propping up the fracked gas industry by manufacturing bloated demand for its dirty product, all with the help of government subsidies and intervention.\\\" \",\n  \"While the IEA report attributes much of the increased demand to a growing number of natural gas power plants in the United States and China, it also points to other factors. U.S. News & World Report outlined the agency's findings:  The industrial sector... also played an outsized role in 2018, with factories, fabricators, and other facilities using gas as both a fuel source and a feedstock to make plastics, fertilizers, and other products\u2014putting industry on track to account for nearly half of global gas consumption by 2024. \",\n  \"The U.S., meanwhile, saw the biggest jump in production last year since 1951, with output soaring by 11.5 percent. That made the U.S. the biggest contributor to gas production growth around the world. IEA executive director Fatih Birol said in a statement announcing the agency's report that \\\"natural gas can contribute to a cleaner global energy system. But it faces its own challenges, including remaining price competitive in emerging markets and reducing methane emissions along the natural gas supply chain.\\\" \",\n  \"Farhana Yamin, a climate attorney and coordinator at the Extinction Rebellion, toldAgence France-Presse, \\\"Given that this polluting fuel can never be 'clean' and is a key driver of climate chaos, the assertion that it can be part of the path to cleaner energy is highly misleading.\\\"  Lorne Stockman, senior research analyst at Oil Change International, also criticized the agency's position on natural gas.  \\\"When it comes to gas, the IEA horse has blinkers on and is heading straight over the cliff of climate disaster,\\\" he told AFP. \\\"Gas is not clean, cheap, or necessary.\\\" \",\n  \"Food & Water Watch's Hauter said, \\\"The IEA's cheerleading of fracked natural gas as some type of global climate solution is foolish and false.\\\" \\\"The time has come to end the madness by ending artificial economic support for the fossil fuel industry, and investing aggressively in truly clean, renewable energy sources like wind and solar,\\\" she added. \\\"The future of our planet depends on it.\\\"\",\n];\n\nconst body = context.document.body;\n\n// The document starts with exactly one paragraph holding the old text;\n// overwrite it in place with the new first paragraph's text.\nconst firstParagraph = body.paragraphs.getFirst();\nfirstParagraph.insertText(paragraphTexts[0], Word.InsertLocation.replace);\n\n// Append the rest of the new paragraphs, one per item, right after the first.\nlet anchor = firstParagraph;\nfor (let i = 1; i < paragraphTexts.length; i++) {\n  anchor = anchor.insertParagraph(paragraphTexts[i], Word.InsertLocation.after);\n}\n\n// Trailing empty paragraph present at the end of the target document.\nanchor.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Replaces the single long paragraph with the new shorter opening paragraph,\n# then adds the remaining article paragraphs (each its own paragraph),\n# followed by one trailing empty paragraph, matching the target diff.\n\n$paragraphTexts = @(\n  'After a frightening study from last week showed that industrial methane emissions have been \"vastly underestimated,\" a new projection Friday that the United States is on track to become the world''s leading exporter of liquefied natural gas within five years provoked warnings that the American fracking boom could \"end hope for climate stability.\" ',\n  'Liquefied natural gas (LNG) is primarily composed of methane, a greenhouse gas that is 84 times more potent than carbon dioxide over a 20-year period. Methane emissions, by some estimates, are responsible for about a quarter of human-caused global warming.  \"Science confirms that gas is a climate killer,\" Wenonah Hauter, executive director of the U.S. advocacy group Food & Water Watch, said in a statement Friday, citing methane''s planet-warming potential. ',\n  'Hauter''s statement came in response to an International Energy Agency (IEA) annual report, released Friday, that featured the new projection about U.S. LNG exports. The IEA report states that global demand for natural gas grew last year at the fastest rate in nearly a decade and is expected to keep growing, \"driven by strong consumption in fast-growing Asian economies and supported by the continued development of the international gas trade.\"  The IEA''s new release came just two days after Food & Water Watch published a report which, as Hauter put it, \"shows that the power, petrochemical, and LNG export industries are propping up the fracked gas industry by manufacturing bloated demand for its dirty product, all with the help of government subsidies and intervention.\" ',\n  'While the IEA report attributes much of the increased demand to a growing number of natural gas power plants in the United States and China, it also points to other factors. U.S. News & World Report outlined the agency''s findings:  The industrial sector... also played an outsized role in 2018, with factories, fabricators, and other facilities using gas as both a fuel source and a feedstock to make plastics, fertilizers, and other products\u2014putting industry on track to account for nearly half of global gas consumption by 2024. ',\n  'The U.S., meanwhile, saw the biggest jump in production last year since 1951, with output soaring by 11.5 percent. That made the U.S. the biggest contributor to gas production growth around the world. IEA executive director Fatih Birol said in a statement announcing the agency''s report that \"natural gas can contribute to a cleaner global energy system. But it faces its own challenges, including remaining price competitive in emerging markets and reducing methane emissions along the natural gas supply chain.\" ',\n  'Farhana Yamin, a climate attorney and coordinator at the Extinction Rebellion, toldAgence France-Presse, \"Given that this polluting fuel can never be ''clean'' and is a key driver of climate chaos, the assertion that it can be part of the path to cleaner energy is highly misleading.\"  Lorne Stockman, senior research analyst at Oil Change International, also criticized the agency''s position on natural gas.  \"When it comes to gas, the IEA horse has blinkers on and is heading straight over the cliff of climate disaster,\" he told AFP. \"Gas is not clean, cheap, or necessary.\" ',\n  'Food & Water Watch''s Hauter said, \"The IEA''s cheerleading of fracked natural gas as some type of global climate solution is foolish and false.\" \"The time has come to end the madness by ending artificial economic support for the fossil fuel industry, and investing aggressively in truly clean, renewable energy sources like wind and solar,\" she added. \"The future of our planet depends on it.\"',\n)\n\n$d = $word.ActiveDocument\n\n# The document starts with exactly one paragraph holding the old text;\n# overwrite it in place with the new first paragraph's text.\n$p = $d.Paragraphs(1)\n$p.Range.Text = $paragraphTexts[0]\n\n# Append the rest of the new paragraphs, one per item, right after the first.\nfor ($i = 1; $i -lt $paragraphTexts.Length; $i++) {\n  $p.Range.InsertParagraphAfter()\n  $p = $d.Paragraphs($i + 1)\n  $p.Range.Text = $paragraphTexts[$i]\n}\n\n# Trailing empty paragraph present at the end of the target document.\n$p.Range.InsertParagraphAfter()\n\n"}
